# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (e.g. "1.00", "9.00")
# must be forced to Text format first so Excel does not convert them into
# real numeric values (which would lose formatting such as trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

# Row 2
$ws.Range("D2").Value2 = "68.596.96"
$ws.Range("E2").Value2 = "  +2.10%  "

# Row 3
$ws.Range("D3").Value2 = "2.523.00"
$ws.Range("E3").Value2 = "  +1.71%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value2 = "  -0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "592.21"
$ws.Range("E5").Value2 = "  +1.25%  "

# Row 6
Set-TextValue $ws.Range("D6") "176.44"
$ws.Range("E6").Value2 = "  +1.82%  "

# Row 7
$ws.Range("E7").Value2 = "  +0.01%  "

# Row 8
$ws.Range("E8").Value2 = "  +0.58%  "

# Row 9
$ws.Range("E9").Value2 = "  +4.40%  "

# Row 10
$ws.Range("E10").Value2 = "  -1.40%  "

# Row 11
Set-TextValue $ws.Range("D11") "4.99"
$ws.Range("E11").Value2 = "  +0.67%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.337"
$ws.Range("E12").Value2 = "  +1.10%  "

# Row 13
$ws.Range("D13").Value2 = "2.939.85"
$ws.Range("E13").Value2 = "  +0.18%  "

# Row 14
Set-TextValue $ws.Range("D14") "25.82"
$ws.Range("E14").Value2 = "  +1.20%  "

# Row 15
$ws.Range("D15").Value2 = "68.331.02"
$ws.Range("E15").Value2 = "  +1.84%  "

# Row 16
$ws.Range("E16").Value2 = "  +0.57%  "

# Row 17
$ws.Range("D17").Value2 = "2.511.69"
$ws.Range("E17").Value2 = "  +2.97%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.02"
$ws.Range("E18").Value2 = "  +0.67%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.49"
$ws.Range("E19").Value2 = "  -0.99%  "

# Row 20
Set-TextValue $ws.Range("D20") "352.08"
$ws.Range("E20").Value2 = "  +0.35%  "

# Row 21
Set-TextValue $ws.Range("D21") "4.15"
$ws.Range("E21").Value2 = "  +3.01%  "

# Row 22
Set-TextValue $ws.Range("D22") "71.37"
$ws.Range("E22").Value2 = "  +3.67%  "

# Row 23
$ws.Range("E23").Value2 = "  -0.02%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.24"
$ws.Range("E24").Value2 = "  +0.32%  "

# Row 25
$ws.Range("E25").Value2 = "  -5.43%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.00"
$ws.Range("E26").Value2 = "  -1.94%  "

# Row 27
$ws.Range("D27").Value2 = "2.588.67"
$ws.Range("E27").Value2 = "  -0.74%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.996"
$ws.Range("E28").Value2 = "  -0.27%  "

# Row 29
$ws.Range("D29").Value2 = "0.0₃0904"
$ws.Range("E29").Value2 = "  -0.66%  "

# Row 30
Set-TextValue $ws.Range("D30") "506.70"
$ws.Range("E30").Value2 = "  +0.66%  "

# Row 31
$ws.Range("E31").Value2 = "  +1.52%  "

# Row 32
$ws.Range("E32").Value2 = "  +2.59%  "

# Row 33
$ws.Range("E33").Value2 = "  +0.87%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.00"
$ws.Range("E34").Value2 = "  +0.02%  "

# Row 35 - B35='Kaspa'
$ws.Range("B35").Value2 = "Kaspa"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D35") "0.121"
$ws.Range("E35").Value2 = "  +2.25%  "

# Row 36 - B36='Monero'
$ws.Range("B36").Value2 = "Monero"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D36") "162.63"
$ws.Range("E36").Value2 = "  -0.15%  "

# Row 37
Set-TextValue $ws.Range("D37") "18.68"
$ws.Range("E37").Value2 = "  -0.11%  "

# Row 38
Set-TextValue $ws.Range("D38") "18.41"
$ws.Range("E38").Value2 = "  +1.22%  "

# Row 39
$ws.Range("E39").Value2 = "  -0.03%  "

# Row 40 - B40='Stacks'
$ws.Range("B40").Value2 = "Stacks"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "1.76"
$ws.Range("E40").Value2 = "  +4.13%  "

# Row 41 - B41='USDe'
$ws.Range("B41").Value2 = "USDe"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value2 = "  +0.04%  "

# Row 42
$ws.Range("E42").Value2 = "  +0.25%  "

# Row 43
Set-TextValue $ws.Range("D43") "4.86"
$ws.Range("E43").Value2 = "  +0.82%  "

# Row 44
$ws.Range("E44").Value2 = "  +1.55%  "

# Row 45
Set-TextValue $ws.Range("D45") "149.77"
$ws.Range("E45").Value2 = "  +4.85%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.56"
$ws.Range("E46").Value2 = "  +2.28%  "

# Row 47 - B47='BabyDogeCoin'
$ws.Range("B47").Value2 = "BabyDogeCoin"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value2 = "0.0₆0260"
$ws.Range("E47").Value2 = "  -0.42%  "

# Row 48 - B48='ARBITRUM'
$ws.Range("B48").Value2 = "ARBITRUM"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D48") "0.520"
$ws.Range("E48").Value2 = "  +1.08%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0739"
$ws.Range("E49").Value2 = "  -0.09%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.59"
$ws.Range("E50").Value2 = "  +0.30%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.582"
$ws.Range("E51").Value2 = "  -0.35%  "

